# Edit script: apply MOSIP Partner Management Requirements content update
# (restores content to the pre-revert state: expands the table from 18 to 20
# data rows, rewrites Module/Feature/Acceptance Criteria/Comments text for most
# rows, and tweaks a couple of cell/row formats.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the table from 18 to 20 data rows (2 new rows appended) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# --- Title row ---
$ws.Range("B1").Value2 = "Partner Management (MISP and E-KYC/Auth Partners)"

# --- Header row (unchanged content, row 2) ---
$ws.Range("B2").Value2 = "Sr No."
$ws.Range("C2").Value2 = "Module"
$ws.Range("D2").Value2 = "Feature"
$ws.Range("E2").Value2 = "Acceptance Criteria"
$ws.Range("F2").Value2 = "Comments"

# --- Data rows 3-20 (set every cell explicitly; clear Comments cells that
#     no longer apply once rows shift) ---
# Row 3
$ws.Range("B3").Value2 = 1
$ws.Range("C3").Value2 = "Kernel"
$ws.Range("D3").Value2 = "MISP ID Generation"
$ws.Range("E3").Value2 = "1. Generate MISP ID as per below logic`na. MISP ID should be of 3 digits (Configurable)`nb. MISP ID should be generated sequentially`nc. MISP ID should be generated incrementally for every request"
$ws.Range("F3").Value2 = "Component already exist as TSP ID generator"

# Row 4
$ws.Range("B4").Value2 = 2
$ws.Range("C4").Value2 = "Kernel"
$ws.Range("D4").Value2 = "MISP License Key Generation"
$ws.Range("E4").Value2 = "1. Generate a License Key as per below logic`na. License Key generation to follow random pattern`nb. License Key should be alphanumeric`nc. Length should be 8 digits (Configurable)`nd. Should be mapped to an expiry"
$ws.Range("F4").ClearContents() | Out-Null

# Row 5
$ws.Range("B5").Value2 = 3
$ws.Range("C5").Value2 = "Kernel"
$ws.Range("D5").Value2 = "MISP License Key Pattern Validation"
$ws.Range("E5").Value2 = "1. Validate length of a License Key as configured and respond as mentioned below`na. If found valid, respond with `"Valid`"`nb. if found invalid, respond with `"Invalid`""
$ws.Range("F5").ClearContents() | Out-Null

# Row 6
$ws.Range("B6").Value2 = 4
$ws.Range("C6").Value2 = "Admin"
$ws.Range("D6").Value2 = "MSIP License Key Expiry Validation"
$ws.Range("E6").Value2 = "1. Validate status of Lisence Key and respond as mentioned below`na. If found expired, respond with `"Your License Key is EXPIRED. Please regenrate a new License Key`"`nb. If found temporarily sespended, respond with `"Your License Key is temporarily SUSPENDED. Please contact MOSIP Administration`"`nc. If found permanently blocked, respond with `"Your License Key is BLOCKED. Please contact MOSIP Administration`""
$ws.Range("F6").ClearContents() | Out-Null

# Row 7
$ws.Range("B7").Value2 = 5
$ws.Range("C7").Value2 = "Admin"
$ws.Range("D7").Value2 = "MISP Registration"
$ws.Range("E7").Value2 = "1. Receive request to register a MISP with follwing parameters`na. MISP Name`nb. MISP Contact Name`nc. MISP Phone`nd. MISP Email ID`n2. Issue and Map MISP ID`n3. Issue and Map Lisence Key`n4. Store the MISP in MOSIP"
$ws.Range("F7").ClearContents() | Out-Null

# Row 8
$ws.Range("B8").Value2 = 6
$ws.Range("C8").Value2 = "Kernel"
$ws.Range("D8").Value2 = "Partner ID Generation"
$ws.Range("E8").Value2 = "1. Generate Partner ID as per below logic`na. Partner ID should be of 4 digits (Configurable)`nb. Partner ID should be generated sequentially`nc. Partner ID should be generated incrementally for every request"
$ws.Range("F8").ClearContents() | Out-Null

# Row 9
$ws.Range("B9").Value2 = 7
$ws.Range("C9").Value2 = "Kernel"
$ws.Range("D9").Value2 = "Partner ID Validation"
$ws.Range("E9").Value2 = "1. Validate length of a Partner ID as configured and respond as mentioned below`na. If found valid, respond with `"Valid`"`nb. if found invalid, respond with `"Invalid`""
$ws.Range("F9").ClearContents() | Out-Null

# Row 10
$ws.Range("B10").Value2 = 8
$ws.Range("C10").Value2 = "Kernel "
$ws.Range("D10").Value2 = "Policy ID Generation"
$ws.Range("E10").Value2 = "1. Generate Policy ID for following policies`na. OTP Trigger `nb. OTP Authentication`nc. Demo Authentication `nd. Biometric Authentication - FMR Data Match `ne. Biometric Authentication - IIR Data Match  `nf. Biometric Authentication - FID Data Match `ng. Static Pin Authentication`nh. eKYC - all combinations of eKYC demo fields `ni. Masked UIN`nj. UIN`n2. Generate Policy id as per below logic`na. Random ID generation`nb. Length should be 10 Digits (Configurable)"
$ws.Range("F10").ClearContents() | Out-Null

# Row 11
$ws.Range("B11").Value2 = 9
$ws.Range("C11").Value2 = "Kernel"
$ws.Range("D11").Value2 = "Policy ID Validation"
$ws.Range("E11").Value2 = "1. Validate length of a Policy ID as configured and respond as mentioned below`na. If found valid, respond with `"Valid`"`nb. if found invalid, respond with `"Invalid`""
$ws.Range("F11").ClearContents() | Out-Null

# Row 12
$ws.Range("B12").Value2 = 10
$ws.Range("C12").Value2 = "Admin"
$ws.Range("D12").Value2 = "Policy ID"
$ws.Range("E12").Value2 = "1. Receive request to retreive policies based on Partner ID and Policy ID`n2. Respond appropirately if Partner ID or Policy ID does not exist"
$ws.Range("F12").ClearContents() | Out-Null

# Row 13
$ws.Range("B13").Value2 = 11
$ws.Range("C13").Value2 = "Admin"
$ws.Range("D13").Value2 = "Partner Registration"
$ws.Range("E13").Value2 = "1. Receive request to register a Partner with follwing parameters`na. Partner Name`nb. Partner Contact Name`nc. Partner Phone`nd. Partner Email ID`n2. Issue and Map Partner ID`n3. Map Policy ID to the Partner`na. Multiple Policies can be mapped to a Partner`nb. A Partner can have a policy for both Auth and E-KYC`n4. Store the Partner in MOSIP"
$ws.Range("F13").ClearContents() | Out-Null

# Row 14
$ws.Range("B14").Value2 = 12
$ws.Range("C14").Value2 = "Admin"
$ws.Range("D14").Value2 = "MISP - Partner Mapping"
$ws.Range("E14").Value2 = "1. Receive a request to map MISP to a Partner with MISP ID and Partner ID as Input`n2. There can ve a many-to-mapping between MISPs and Partners"
$ws.Range("F14").ClearContents() | Out-Null

# Row 15
$ws.Range("B15").Value2 = 13
$ws.Range("C15").Value2 = "Admin"
$ws.Range("D15").Value2 = "Partner Certiicate Validation"
$ws.Range("E15").Value2 = "1. Receive certificate from Partner`n2. Verify CA Authority of the certificate"
$ws.Range("F15").ClearContents() | Out-Null

# Row 16
$ws.Range("B16").Value2 = 14
$ws.Range("C16").Value2 = "Admin"
$ws.Range("D16").Value2 = "Partner Certificate Signing and RE-issueing"
$ws.Range("E16").Value2 = "1. Receive certificate from Partner during Partner Registration`n2. Sign the Partner Certificate with MOSIP Private Key and issue a certificate chain`n3. Re-issue certficate back to the Partner`n4. Private key to change priodically as per the Key Rotation Policy set by admin"
$ws.Range("F16").ClearContents() | Out-Null

# Row 17
$ws.Range("B17").Value2 = 15
$ws.Range("C17").Value2 = "Admin"
$ws.Range("D17").Value2 = "Distribution of Public Key to Partners"
$ws.Range("E17").Value2 = "1. Distribute Public Key to Partners correspinding to the Private Key used to signed the Certificate`n2. Public key needs to be distributed priodically whenever the Private Key is rotated"
$ws.Range("F17").ClearContents() | Out-Null

# Row 18
$ws.Range("B18").Value2 = 16
$ws.Range("C18").Value2 = "Admin"
$ws.Range("D18").Value2 = "Device Registration"
$ws.Range("E18").Value2 = "TBD"
$ws.Range("F18").Value2 = "Yet to analyzed"

# Row 19
$ws.Range("B19").Value2 = 17
$ws.Range("C19").Value2 = "Admin"
$ws.Range("D19").Value2 = "Device Provider Registration"
$ws.Range("E19").Value2 = "TBD"
$ws.Range("F19").Value2 = "Yet to analyzed"

# Row 20
$ws.Range("B20").Value2 = 18
$ws.Range("C20").Value2 = "Admin"
$ws.Range("D20").Value2 = "RD Service Registration"
$ws.Range("E20").Value2 = "TBD"
$ws.Range("F20").Value2 = "Yet to analyzed"

# --- Row heights (to match autofit-wrapped content) ---
$ws.Rows(1).RowHeight = 21.5
$ws.Rows(2).RowHeight = 15
$ws.Rows(3).RowHeight = 58
$ws.Rows(4).RowHeight = 72.5
$ws.Rows(5).RowHeight = 58
$ws.Rows(6).RowHeight = 101.5
$ws.Rows(7).RowHeight = 116
$ws.Rows(8).RowHeight = 58
$ws.Rows(9).RowHeight = 58
$ws.Rows(10).RowHeight = 203
$ws.Rows(11).RowHeight = 58
$ws.Rows(12).RowHeight = 43.5
$ws.Rows(13).RowHeight = 145
$ws.Rows(14).RowHeight = 43.5
$ws.Rows(15).RowHeight = 29
$ws.Rows(16).RowHeight = 87
$ws.Rows(17).RowHeight = 58
$ws.Rows(18).RowHeight = 14.5
$ws.Rows(19).RowHeight = 14.5
$ws.Rows(20).RowHeight = 14.5

# --- Special formatting: Acceptance Criteria cell on the new "Partner Certificate
# Signing and RE-issueing" row gets left-aligned wrapped text (distinct from the
# default wrap-only style used by the rest of the Acceptance Criteria column) ---
$ws.Range("E16").HorizontalAlignment = -4131

# --- Scroll back to the top and restore the active selection ---
$ws.Range("E7").Select() | Out-Null